# Automatische test-sync: 2025-06-19 15:00:10
$wb = $excel.ActiveWorkbook

# "Logs" sheet holds the raw mail log rows
$logs = $wb.Worksheets.Item("Logs")

# Append the new incoming mail as row 15
$logs.Cells.Item(15, 1).Value = "Vragen over samenwerking"
$logs.Cells.Item(15, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(15, 3).Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Cells.Item(15, 4).Value = "Overig"
$logs.Cells.Item(15, 6).Value = "2025-06-19 14:58:10"
$logs.Cells.Item(15, 7).Value = "Nee"

# Extend the conditional-formatting ranges to cover the new row
$catRules = $logs.Range("D2:D14").FormatConditions
$catRules.Item(1).ModifyAppliesToRange($logs.Range("D2:D15"))

$answeredRules = $logs.Range("G2:G14").FormatConditions
$answeredRules.Item(1).ModifyAppliesToRange($logs.Range("G2:G15"))

# "Dashboard" sheet keeps a per-category tally; bump the "Overig" count
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 8
